$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 6.857917524323606
$ws.Range("E3").Value = 6.503119830300792
$ws.Range("E4").Value = 6.218080370165573
$ws.Range("E5").Value = 6.341921423272087
$ws.Range("E6").Value = 5.968444520820012
$ws.Range("E7").Value = 5.776644732809495
$ws.Range("E8").Value = 4.556109122745393
$ws.Range("E9").Value = 3.66639353244369
$ws.Range("E10").Value = 3.082808567070776
$ws.Range("E11").Value = 3.091707788714488
$ws.Range("E12").Value = 3.039115502164408
$ws.Range("E13").Value = 3.11801755371366
$ws.Range("E14").Value = 3.208169496592049
$ws.Range("E15").Value = 3.072004752223021
$ws.Range("E16").Value = 3.42123495093697
$ws.Range("E17").Value = 3.418196937883507
$ws.Range("E18").Value = 3.326729283189561
$ws.Range("E19").Value = 3.383011347750638
$ws.Range("E20").Value = 3.195954794956731
$ws.Range("E21").Value = 2.752190229701109
$ws.Range("E22").Value = 2.939456748273071
$ws.Range("E23").Value = 2.867299028563491
$ws.Range("E24").Value = 3.140089936237962
$ws.Range("E25").Value = 3.092905690763107
$ws.Range("E26").Value = 2.959371991001731
$ws.Range("E27").Value = 3.155688164795017
$ws.Range("E28").Value = 3.045540524691801
$ws.Range("E29").Value = 3.413056921577232
$ws.Range("E30").Value = 3.602665687103285
$ws.Range("E31").Value = 4.039625292044262
$ws.Range("E32").Value = 4.667291911210431
$ws.Range("E33").Value = 5.467678370038672
$ws.Range("E34").Value = 5.515613192695801
$ws.Range("E35").Value = 4.493179863750797
$ws.Range("E36").Value = 4.509818410644426
